$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "Dates" column (column E) entirely; all columns to the right shift left.
$ws.Range("E1").EntireColumn.Delete()
